$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 69
$ws.Range("H69").Value = 4873.421
$ws.Range("I69").Value = 4531.6665
$ws.Range("J69").Value = 5031.154
$ws.Range("K69").Value = 13594.9995
$ws.Range("L69").Value = 15093.462
$ws.Range("M69").Value = -12720.9995
$ws.Range("N69").Value = -16841.462

# Row 72
$ws.Range("H72").Value = 4873.421
$ws.Range("I72").Value = 4531.6665
$ws.Range("J72").Value = 5031.154
$ws.Range("K72").Value = 40784.9985
$ws.Range("L72").Value = 45280.38600000001
$ws.Range("M72").Value = -36416.9985
$ws.Range("N72").Value = -54016.38600000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4360.8
$ws.Range("I32").Value = 2948.2144
$ws.Range("J32").Value = 10011.143
$ws.Range("K32").Value = 2948.2144
$ws.Range("L32").Value = 10011.143
$ws.Range("M32").Value = -2661.2144
$ws.Range("N32").Value = -10585.143

# Row 63
$ws.Range("H63").Value = 166668580
$ws.Range("I63").Value = 166668580
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 166668580
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -166667894
$ws.Range("N63").ClearContents()

# Row 66
$ws.Range("H66").Value = 166668580
$ws.Range("I66").Value = 166668580
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 833342900
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -833339468
$ws.Range("N66").ClearContents()

# Row 88
$ws.Range("H88").Value = 2900
$ws.Range("I88").Value = 2900
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 2900
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -2494
$ws.Range("N88").ClearContents()

# Row 91
$ws.Range("H91").Value = 2900
$ws.Range("I91").Value = 2900
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 2900
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -1496
$ws.Range("N91").ClearContents()

# Row 102
$ws.Range("H102").Value = 3087963
$ws.Range("I102").Value = 3705104.5
$ws.Range("K102").Value = 3705104.5
$ws.Range("M102").Value = -3703482.5

$ws = $wb.Worksheets.Item("BSM")
# Row 54
$ws.Range("H54").Value = 17750
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 17750
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 17750
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -18718

# Row 86
$ws.Range("H86").Value = 22223942
$ws.Range("I86").Value = 30304566
$ws.Range("J86").Value = 2226.75
$ws.Range("K86").Value = 30304566
$ws.Range("L86").Value = 2226.75
$ws.Range("M86").Value = -30303443
$ws.Range("N86").Value = -4472.75

# Row 89
$ws.Range("H89").Value = 22223942
$ws.Range("I89").Value = 30304566
$ws.Range("J89").Value = 2226.75
$ws.Range("K89").Value = 151522830
$ws.Range("L89").Value = 11133.75
$ws.Range("M89").Value = -151517214
$ws.Range("N89").Value = -22365.75

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 162.4
$ws.Range("I7").Value = 153.09091
$ws.Range("J7").Value = 188
$ws.Range("K7").Value = 153.09091
$ws.Range("L7").Value = 188
$ws.Range("M7").Value = -40.09091000000001
$ws.Range("N7").Value = -414

# Row 57
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

# Row 99
$ws.Range("H99").Value = 12504500
$ws.Range("I99").Value = 1866.1666
$ws.Range("K99").Value = 1866.1666
$ws.Range("M99").Value = -368.1666

# Row 126
$ws.Range("H126").Value = 12504500
$ws.Range("I126").Value = 1866.1666
$ws.Range("K126").Value = 5598.4998
$ws.Range("M126").Value = -3128.4998

# Row 141
$ws.Range("H141").Value = 21746
$ws.Range("J141").Value = 21746
$ws.Range("L141").Value = 21746
$ws.Range("N141").Value = -32106

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 286654
$ws.Range("I5").Value = 668.4286
$ws.Range("J5").Value = 429646.78
$ws.Range("K5").Value = 2005.2858
$ws.Range("L5").Value = 1288940.34
$ws.Range("M5").Value = -1893.2858
$ws.Range("N5").Value = -1289164.34

# Row 8
$ws.Range("H8").Value = 714378.4399999999
$ws.Range("I8").Value = 714378.4399999999
$ws.Range("K8").Value = 2143135.32
$ws.Range("M8").Value = -2142996.32

# Row 38
$ws.Range("H38").Value = 3846521.5
$ws.Range("I38").Value = 6667015
$ws.Range("J38").Value = 394.18182
$ws.Range("K38").Value = 20001045
$ws.Range("L38").Value = 1182.54546
$ws.Range("M38").Value = -20000698
$ws.Range("N38").Value = -1876.54546

# Row 96
$ws.Range("H96").Value = 4700
$ws.Range("J96").Value = 4700
$ws.Range("L96").Value = 14100
$ws.Range("N96").Value = -18218

# Row 97
$ws.Range("H97").Value = 14286387
$ws.Range("I97").Value = 50000250
$ws.Range("J97").Value = 841.6
$ws.Range("K97").Value = 150000750
$ws.Range("L97").Value = 2524.8
$ws.Range("M97").Value = -150000254
$ws.Range("N97").Value = -3516.8

# Row 122
$ws.Range("H122").Value = 763.55554
$ws.Range("I122").Value = 611.1539
$ws.Range("J122").Value = 1159.8
$ws.Range("K122").Value = 5500.3851
$ws.Range("L122").Value = 10438.2
$ws.Range("M122").Value = -3050.3851
$ws.Range("N122").Value = -15338.2

# Row 135
$ws.Range("H135").Value = 286654
$ws.Range("I135").Value = 668.4286
$ws.Range("J135").Value = 429646.78
$ws.Range("K135").Value = 6015.8574
$ws.Range("L135").Value = 3866821.02
$ws.Range("M135").Value = -3480.8574
$ws.Range("N135").Value = -3871891.02

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 16667677
$ws.Range("I46").Value = 27778534
$ws.Range("J46").Value = 1392.375
$ws.Range("K46").Value = 27778534
$ws.Range("L46").Value = 1392.375
$ws.Range("M46").Value = -27778346
$ws.Range("N46").Value = -1768.375

# Row 56
$ws.Range("H56").Value = 13616.667
$ws.Range("J56").Value = 13616.667
$ws.Range("L56").Value = 13616.667
$ws.Range("N56").Value = -14998.667

# Row 122
$ws.Range("H122").Value = 7939183.5
$ws.Range("I122").Value = 8931269
$ws.Range("K122").Value = 26793807
$ws.Range("M122").Value = -26791357

$ws = $wb.Worksheets.Item("WVR")
# Row 58
$ws.Range("H58").Value = 15000
$ws.Range("J58").Value = 15000
$ws.Range("L58").Value = 15000
$ws.Range("N58").Value = -15616

# Row 81
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()

# Row 136
$ws.Range("H136").Value = 2119
$ws.Range("I136").Value = 2228.5334
$ws.Range("J136").Value = 1790.4
$ws.Range("K136").Value = 6685.600199999999
$ws.Range("L136").Value = 5371.200000000001
$ws.Range("M136").Value = -4135.600199999999
$ws.Range("N136").Value = -10471.2
